$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.389.68"
$ws.Cells.Item(2, 5).Value = "  +4.16%  "
$ws.Cells.Item(3, 4).Value = "1.586.36"
$ws.Cells.Item(3, 5).Value = "  +1.25%  "
$ws.Cells.Item(4, 5).Value = "  -0.34%  "
$ws.Cells.Item(5, 4).Value = "'214.35"
$ws.Cells.Item(5, 5).Value = "  +1.68%  "
$ws.Cells.Item(6, 5).Value = "  +1.00%  "
$ws.Cells.Item(7, 5).Value = "  -0.37%  "
$ws.Cells.Item(8, 4).Value = "'23.81"
$ws.Cells.Item(8, 5).Value = "  +7.78%  "
$ws.Cells.Item(9, 5).Value = "  +0.83%  "
$ws.Cells.Item(10, 4).Value = "'0.0600"
$ws.Cells.Item(10, 5).Value = "  +0.44%  "
$ws.Cells.Item(11, 4).Value = "'0.0890"
$ws.Cells.Item(11, 5).Value = "  +2.34%  "
$ws.Cells.Item(12, 4).Value = "1.816.81"
$ws.Cells.Item(12, 5).Value = "  +1.47%  "
$ws.Cells.Item(13, 4).Value = "1.587.27"
$ws.Cells.Item(13, 5).Value = "  +1.09%  "
$ws.Cells.Item(14, 5).Value = "  +0.21%  "
$ws.Cells.Item(15, 5).Value = "  +2.44%  "
$ws.Cells.Item(16, 4).Value = "28.378.16"
$ws.Cells.Item(16, 5).Value = "  +4.21%  "
$ws.Cells.Item(17, 4).Value = "'63.39"
$ws.Cells.Item(17, 5).Value = "  +1.95%  "
$ws.Cells.Item(18, 4).Value = "'231.71"
$ws.Cells.Item(18, 5).Value = "  +6.62%  "
$ws.Cells.Item(19, 5).Value = "  +0.93%  "
$ws.Cells.Item(20, 4).Value = "'7.48"
$ws.Cells.Item(20, 5).Value = "  -0.32%  "
$ws.Cells.Item(21, 5).Value = "  -0.18%  "
$ws.Cells.Item(22, 4).Value = "'4.12"
$ws.Cells.Item(22, 5).Value = "  -0.73%  "
$ws.Cells.Item(23, 4).Value = "'9.40"
$ws.Cells.Item(23, 5).Value = "  +1.88%  "
$ws.Cells.Item(24, 5).Value = "  +1.05%  "
$ws.Cells.Item(25, 4).Value = "'152.44"
$ws.Cells.Item(25, 5).Value = "  -0.62%  "
$ws.Cells.Item(26, 4).Value = "'15.25"
$ws.Cells.Item(26, 5).Value = "  +1.29%  "
$ws.Cells.Item(27, 4).Value = "'6.61"
$ws.Cells.Item(27, 5).Value = "  -0.32%  "
$ws.Cells.Item(28, 4).Value = "'0.107"
$ws.Cells.Item(28, 5).Value = "  +0.89%  "
$ws.Cells.Item(29, 5).Value = "  -0.29%  "
$ws.Cells.Item(30, 5).Value = "  +0.48%  "
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 5).Value = "  -0.18%  "
$ws.Cells.Item(33, 5).Value = "  -0.46%  "
$ws.Cells.Item(34, 4).Value = "1.408.60"
$ws.Cells.Item(34, 5).Value = "  -2.43%  "
$ws.Cells.Item(35, 5).Value = "  -1.24%  "
$ws.Cells.Item(36, 4).Value = "'1.05"
$ws.Cells.Item(36, 5).Value = "  -4.54%  "
$ws.Cells.Item(37, 4).Value = "'2.34"
$ws.Cells.Item(37, 5).Value = "  +0.06%  "
$ws.Cells.Item(38, 5).Value = "  +0.44%  "
$ws.Cells.Item(39, 5).Value = "  +8.55%  "
$ws.Cells.Item(40, 5).Value = "  +1.60%  "
$ws.Cells.Item(41, 4).Value = "'0.816"
$ws.Cells.Item(41, 5).Value = "  +0.83%  "
$ws.Cells.Item(42, 4).Value = "'5.76"
$ws.Cells.Item(42, 5).Value = "  -2.07%  "
$ws.Cells.Item(43, 5).Value = "  -0.29%  "
$ws.Cells.Item(44, 4).Value = "'0.983"
$ws.Cells.Item(44, 5).Value = "  -1.93%  "
$ws.Cells.Item(45, 4).Value = "'1.83"
$ws.Cells.Item(45, 5).Value = "  +5.72%  "
$ws.Cells.Item(46, 4).Value = "'64.54"
$ws.Cells.Item(46, 5).Value = "  +0.06%  "
$ws.Cells.Item(47, 4).Value = "1.728.38"
$ws.Cells.Item(47, 5).Value = "  +1.42%  "
$ws.Cells.Item(48, 4).Value = "'87.60"
$ws.Cells.Item(48, 5).Value = "  +1.88%  "
$ws.Cells.Item(49, 5).Value = "  +5.39%  "
$ws.Cells.Item(50, 5).Value = "  -0.66%  "
$ws.Cells.Item(51, 4).Value = "'39.31"
$ws.Cells.Item(51, 5).Value = "  +16.04%  "
